$d = $word.ActiveDocument

# 1) "recorded_at" row, Format column: "timestamp with time zone" -> "date-time"
#    (leaves the sibling "Read only" run untouched)
$d.Content.Find.Execute(
    "timestamp with time zone", $true, $false, $false, $false, $false,
    $true, 1, $false, "date-time", 2
)

# 2) First table's column widths were re-flowed (e.g. Word recalculated the
#    grid after the shorter cell text above). Word.Column.Width is expressed
#    in points, while the underlying <w:gridCol w:w="..."/> is in twips
#    (1 pt = 20 twips), so divide the target twip values by 20.
$t = $d.Tables.Item(1)
$t.Columns.Item(1).Width = 1762 / 20
$t.Columns.Item(2).Width = 4369 / 20
$t.Columns.Item(3).Width = 968 / 20
$t.Columns.Item(4).Width = 819 / 20
